# Generate Report for Handback
#
# Marks the zh-cn and de-de localization rows as handed back: fills in the
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns, updates the Status text from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it is shown (Overview summary
# + both language sheets), widens a few columns that now hold longer text,
# and adds a hyperlink on the new "Latest Target File" cell pointing at the
# same source doc as the existing "Source File Name" hyperlink.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$docName = "127c6d9f-82d3-41eb-88e8-925348b722ca.md"
$docUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/f055449e56b44a9b9cd9b0b9b81deb8dcc4599e2/e2e/127c6d9f-82d3-41eb-88e8-925348b722ca.md"

$handedBackStatus = "Handed back: in sync with en-US"

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# Update every cell that shows the status (Overview!E2/F2 summary + the
# per-language Status column) so the shared string collapses onto the new
# text everywhere it is used.
$wsOverview.Range("E2").Value = $handedBackStatus
$wsOverview.Range("F2").Value = $handedBackStatus
$wsZhCn.Range("C2").Value = $handedBackStatus
$wsDeDe.Range("C2").Value = $handedBackStatus

# --- zh-cn row: fill in target/handback file + datetime ---
$wsZhCn.Range("J2").Value = "127c6d9f-82d3-41eb-88e8-925348b722ca.ade85a9222efaa916ffcbb730b3a1dffce132907.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-13 23:17:19"
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(2, 9), $docUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $docName)

# --- de-de row: fill in target/handback file + datetime ---
$wsDeDe.Range("J2").Value = "127c6d9f-82d3-41eb-88e8-925348b722ca.ade85a9222efaa916ffcbb730b3a1dffce132907.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-13 23:17:29"
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(2, 9), $docUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $docName)

# --- Column widths: widen columns that now hold the longer handback text ---
# The headless ColumnWidth setter pads whatever is assigned by a fixed
# ~5/6-character offset before it is written back out to the sheet's <col>
# element, so the literal values below are pre-compensated to land on the
# intended widths (29.9777047293527 / 40) once that padding is re-applied.
$widthOffset = 0.8333333333335986
$wideColWidth = 29.9777047293527 - $widthOffset
$fullColWidth = 40.0 - $widthOffset

$wsOverview.Columns("E:E").ColumnWidth = $wideColWidth
$wsOverview.Columns("F:F").ColumnWidth = $wideColWidth

$wsZhCn.Columns("C:C").ColumnWidth = $wideColWidth
$wsZhCn.Columns("I:I").ColumnWidth = $fullColWidth
$wsZhCn.Columns("J:J").ColumnWidth = $fullColWidth

$wsDeDe.Columns("C:C").ColumnWidth = $wideColWidth
$wsDeDe.Columns("I:I").ColumnWidth = $fullColWidth
$wsDeDe.Columns("J:J").ColumnWidth = $fullColWidth
